$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: Black theme
$ws.Range("A5").Value = "Black"
$ws.Range("B5").Value = "https://www.dreamstime.com/casino-theme-decorative-design-elements-chalkboard-gambling-symbols-casino-theme-decorative-design-elements-chalkboard-image105692673"
$ws.Range("C5").Value = "Custom color in css"
$ws.Range("B5:C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 29.5

# Row 6: Glow theme
$ws.Range("A6").Value = "Glow"
$ws.Range("B6").Value = "https://www.shutterstock.com/image-illustration/dark-bright-pattern-playing-card-symbols-756655723?src=979Kn-kXZpVNhVXiNDOEqw-7-86"
$ws.Range("C6").Value = "Custom color in css"
$ws.Range("B6:C6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 27

$ws.Range("C6").Select()
